$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to the batter's name
$ws.Name = "Moises Henriques"

# Insert a new column A ("matchNo"), shifting teamName..result from A..L to B..M
$ws.Columns("A").Insert()
$ws.Range("A1").Value = "matchNo"

# Insert a new row above the existing data row (old row 2 data -> becomes row 3)
$ws.Rows("2").Insert()

# ---- Row 2: "53rd" match ----
$ws.Range("A2").Value = "53rd"
$ws.Range("B2").Value = "Punjab Kings"
$ws.Range("C2").Value = "Moises Henriques"
$ws.Range("D2").Value = ""
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "3"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "0"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "100.00"
$ws.Range("J2").Value = "Chennai Super Kings"
$ws.Range("K2").Value = "Dubai (DSC)"
$ws.Range("L2").Value = "October 07"
$ws.Range("M2").Value = "Punjab Kings won by 6 wickets (with 42 balls remaining)"

# ---- Row 3: formerly row 2 ("48th" match) -> just fill the new matchNo cell ----
$ws.Range("A3").Value = "48th"

# ---- Row 4: "21st" match ----
$ws.Range("A4").Value = "21st"
$ws.Range("B4").Value = "Punjab Kings"
$ws.Range("C4").Value = "Moises Henriques"
$ws.Range("D4").Value = "b Narine"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "3"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "0"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "0"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "66.66"
$ws.Range("J4").Value = "Kolkata Knight Riders"
$ws.Range("K4").Value = "Ahmedabad"
$ws.Range("L4").Value = "April 26"
$ws.Range("M4").Value = "KKR won by 5 wickets (with 20 balls remaining)"

# ---- Row 5: "14th" match ----
$dagger = [char]0x2020
$ws.Range("A5").Value = "14th"
$ws.Range("B5").Value = "Punjab Kings"
$ws.Range("C5").Value = "Moises Henriques"
$ws.Range("D5").Value = "st " + $dagger + "Bairstow b Abhishek Sharma"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "14"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "17"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "0"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "0"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "82.35"
$ws.Range("J5").Value = "Sunrisers Hyderabad"
$ws.Range("K5").Value = "Chennai"
$ws.Range("L5").Value = "April 21"
$ws.Range("M5").Value = "Sunrisers won by 9 wickets (with 8 balls remaining)"

Write-Host "Applied edits"
